$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# 1. Intro paragraph: "HTML, PHP, Java, JavaScript," -> "HTML, PHP, JavaScript,"
Replace-Text "HTML, PHP, Java, JavaScript, " "HTML, PHP, JavaScript, "

# 2. Technology paragraph, first Java -> JavaScript (before ", és")
Replace-Text "a következő technológiákra építi: HTML, PHP, Java, és " "a következő technológiákra építi: HTML, PHP, JavaScript, és "

# 3. Technology paragraph, second Java -> JavaScript ("A Java a webalkalmazások")
Replace-Text "a felhasználói adatok megjelenítése. A Java a webalkalmazások" "a felhasználói adatok megjelenítése. A JavaScript a webalkalmazások"

# 4. 11.14 paragraph: remove "marketing, " and "Projekt célok finalizálása: "
Replace-Text "tesztelés, marketing, stb. Projekt célok finalizálása: A csapat közösen" "tesztelés, stb. A csapat közösen"

# 5. 11.15 paragraph: remove "Feladatok delegálása: "
Replace-Text "szétosztása: Feladatok delegálása: A csapat tagjai között" "szétosztása: A csapat tagjai között"

# 6. 11.16 paragraph, first Java -> JavaScript
Replace-Text "kódolás: HTML, PHP, Java, MySQL integráció" "kódolás: HTML, PHP, JavaScript, MySQL integráció"

# 7 & 8. 11.16 paragraph, second Java -> JavaScript and typo fix hozjuk -> hozzuk
Replace-Text "technológiák (HTML, PHP, Java, MySQL) segítségével. HTML: A weboldal alapstruktúráját és dizájnját hozjuk létre" "technológiák (HTML, PHP, JavaScript, MySQL) segítségével. HTML: A weboldal alapstruktúráját és dizájnját hozzuk létre"

# 9. 11.16 paragraph, third Java -> JavaScript
Replace-Text "felhasználói interakciók biztosítására. Java: A felhasználói élmény" "felhasználói interakciók biztosítására. JavaScript: A felhasználói élmény"

# 10. 12.07 paragraph: remove " és finomhangolás" from title
Replace-Text "Első tesztelés és finomhangolás: Első tesztelési fázis" "Első tesztelés: Első tesztelési fázis"

# 11. 12.07 paragraph: remove "Bugs és " before "problémák keresése" (capitalize P)
Replace-Text "javítása érdekében. Bugs és problémák keresése:" "javítása érdekében. Problémák keresése:"

# 12. 12.07 paragraph: "ellenőrzése és finomhangolása is történik." -> "ellenőrzése is megtörténik."
Replace-Text "amelyek ellenőrzése és finomhangolása is történik. " "amelyek ellenőrzése is megtörténik. "

# 14. 12.16 paragraph: remove " és UX finomhangolás: Dizájn finalizálása" and "(UX) "
Replace-Text "Weboldal dizájn és UX finomhangolás: Dizájn finalizálása: A felhasználói élmény (UX) és a weboldal" "Weboldal dizájn: A felhasználói élmény  és a weboldal"

# 13. Move the _GoBack bookmark from the empty paragraph (after 12.07) into the
#     12.16 paragraph, positioned right after "A felhasználói élmény " (between the
#     two spaces left behind by the removed "(UX)").
$anchor = $d.Content.Find.Execute("A felhasználói élmény  és a weboldal vizuális megjelenése", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$marker = $d.Content
$found = $marker.Find.Execute("A felhasználói élmény  és a weboldal vizuális megjelenése (UI) finomhangolása", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $marker.Start + [string]"A felhasználói élmény ".Length
$bmRange = $d.Range($start, $start)
$d.Bookmarks.Add("_GoBack", $bmRange)
